$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

$newRow = 81

# Copy the date formatting from the row above so the new date cells reuse
# the existing date style instead of creating a new number format.
$ws.Range("C80").Copy()
$ws.Range("C$newRow").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("D80").Copy()
$ws.Range("D$newRow").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("A$newRow").Value = "The Vaccine Race"
$ws.Range("B$newRow").Value = "Meredith Wadman"
$ws.Range("C$newRow").Value = Get-Date -Year 2020 -Month 5 -Day 22 -Hour 0 -Minute 0 -Second 0
$ws.Range("D$newRow").Value = Get-Date -Year 2020 -Month 5 -Day 28 -Hour 0 -Minute 0 -Second 0
$ws.Range("E$newRow").Value = "vaccines;science;medicine;disease;virus;history"
$ws.Range("F$newRow").Value = "Audio"
$ws.Range("G$newRow").Value = "19 Hours 19 Mins"

$ws.Application.CutCopyMode = $false

$nextCell = $ws.Range("A82")
$ws.Application.ActiveWindow.ScrollRow = 62
$nextCell.Select() | Out-Null
